$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 60 with quarterly data for 01-07-2021.
# Column A holds a date-like label ("01-07-2021") that must be stored as
# literal text (matching the existing Serie column), not auto-converted to
# a date serial number. Temporarily mark the cell as Text before typing the
# value, then clear the formatting back off so the cell ends up with no
# explicit style, same as the rest of the Serie column.
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "01-07-2021"
$ws.Range("A60").ClearFormats()
$ws.Range("B60").Value = 2482
$ws.Range("C60").Value = 2333
$ws.Range("D60").Value = 0
$ws.Range("F60").Value = 148
